# EPBDS-12620 Difference in error response structure between kafka and webservice call
#
# The "Good Night" sample return value on the Rules decision table is
# replaced with the literal text  = error("fail")  so the test data drives
# an error response (used to compare Kafka vs. webservice error payloads).
# A leading apostrophe forces Excel to store it as literal text (quote
# prefix) instead of parsing it as a formula.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$ws.Range("E11").Value = "'= error(""fail"")"
